$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RegisterNewUser")

# Remove all existing hyperlinks on the sheet (Email column hyperlinks plus the
# Password-column hyperlinks that will be re-created after the column shift).
$ws.Cells.Item(4, 5).Hyperlinks.Delete()

# Delete the "Email" column (column D) entirely - this shifts the old
# "Password" column (E) into column D.
$ws.Columns("D").Delete()

# Re-create the single remaining hyperlink (originally on the Email column,
# now on the shifted Password column D4), preserving its mismatched display
# text while keeping the cell's real value as the password string.
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:1234567@", "", "", "1234567@")
$ws.Range("D4").Value = "1234567!"
$ws.Range("D4").ClearFormats()

# Match the post-edit selection state (whole column D selected).
$ws.Columns("D").Select()
